$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style template cells (known good style indices 1/2/3) ---
$styleA = $ws.Cells.Item(3, 1)
$styleB = $ws.Cells.Item(3, 2)
$styleC = $ws.Cells.Item(3, 3)

# Row 1
$ws.Cells.Item(1, 2).Value = "Ementa atual:"
$ws.Cells.Item(1, 3).Value = "Ementa modificada (dados modificados em vermelho):"

# Row 2
$ws.Cells.Item(2, 2).Value = "LOT2041"
$ws.Cells.Item(2, 3).Value = "LOT2041"

# Row 3
$ws.Cells.Item(3, 1).Value = "Nome:"
$ws.Cells.Item(3, 2).Value = " Engenharia Bioquímica"
$ws.Cells.Item(3, 3).Value = " Engenharia Bioquímica"

# Row 4
$ws.Cells.Item(4, 1).Value = "Name:"
$ws.Cells.Item(4, 2).Value = "Biochemical Engineering"
$ws.Cells.Item(4, 3).Value = "Biochemical Engineering"

# Row 5
$ws.Cells.Item(5, 1).Value = "Créditos-aula:"
$ws.Cells.Item(5, 2).Value = "4"
$ws.Cells.Item(5, 3).Value = "4"

# Row 6
$ws.Cells.Item(6, 1).Value = "Créditos-trabalho"
$ws.Cells.Item(6, 2).Value = "0"
$ws.Cells.Item(6, 3).Value = "0"

# Row 7
$ws.Cells.Item(7, 1).Value = "Carga horária:"
$ws.Cells.Item(7, 2).Value = "60 h"
$ws.Cells.Item(7, 3).Value = "60 h"

# Row 8
$ws.Cells.Item(8, 1).Value = "Ativação:"
$ws.Cells.Item(8, 2).Value = "01/01/2019"
$ws.Cells.Item(8, 3).Value = "01/01/2019"

# Row 9
$ws.Cells.Item(9, 1).Value = "Semestre ideal:"
$ws.Cells.Item(9, 2).Value = "EQD-8,EQN-9"
$ws.Cells.Item(9, 3).Value = "EQD-8,EQN-9"

# Row 10
$ws.Cells.Item(10, 1).Value = "Objetivos:"
$ws.Cells.Item(10, 2).Value = "Capacitar o aluno para aplicar os conceitos de Engenharia aos Processos Biológicos e para identificar a relevância dos processos microbianos em escala industrial. Fornecer conhecimentos, técnicas e métodos de base científica ou prática para uma melhor compreensão dos aspectos cinéticos de um bioprocesso em suas diferentes formas de condução (regime descontínuo, descontinuo- alimentado e contínuo), assim como dos conceitos fundamentais para o desenvolvimento da etapa de esterilização de um bioprocesso."
$ws.Cells.Item(10, 3).Value = "Capacitar o aluno para aplicar os conceitos de Engenharia aos Processos Biológicos e para identificar a relevância dos processos microbianos em escala industrial. Fornecer conhecimentos, técnicas e métodos de base científica ou prática para uma melhor compreensão dos aspectos cinéticos de um bioprocesso em suas diferentes formas de condução (regime descontínuo, descontinuo- alimentado e contínuo), assim como dos conceitos fundamentais para o desenvolvimento da etapa de esterilização de um bioprocesso."

# Row 11
$ws.Cells.Item(11, 1).Value = "Objectives:"

# Row 12
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"

# Row 13
$ws.Cells.Item(13, 1).ClearContents()
$ws.Cells.Item(13, 2).Value = "1112574 - Inês Conceição Roberto"
$ws.Cells.Item(13, 3).Value = "1112574 - Inês Conceição Roberto"
$ws.Rows.Item(13).EntireRow.AutoFit()

# Row 14
$ws.Cells.Item(14, 1).ClearContents()
$ws.Cells.Item(14, 2).Value = "1097178 - João Batista de Almeida e Silva"
$ws.Cells.Item(14, 3).Value = "1097178 - João Batista de Almeida e Silva"
$ws.Rows.Item(14).EntireRow.AutoFit()

# Row 15
$ws.Cells.Item(15, 1).Value = "Programa resumido:"
$ws.Cells.Item(15, 2).Value = "Características do material biológico; Cinética de processos fermentativos, Formas decondução dos processos fermentativos, esterilização em bioprocessos."
$ws.Cells.Item(15, 3).Value = "Características do material biológico; Cinética de processos fermentativos, Formas decondução dos processos fermentativos, esterilização em bioprocessos."
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Cells.Item(16, 1).Value = "Short syllabus:"
$ws.Cells.Item(16, 2).Value = "Characteristics of biological material; Kinetics of fermentative processes; Operation modes of fermentative processes; Sterilization in bioprocess."
$ws.Cells.Item(16, 3).Value = "Characteristics of biological material; Kinetics of fermentative processes; Operation modes of fermentative processes; Sterilization in bioprocess."
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Cells.Item(17, 1).Value = "Programa:"
$styleB.Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4122)
$ws.Cells.Item(17, 2).Value = "Características do material biológico – Tipos de microrganismos, morfologia e estrutura celular, nutrição e crescimento microbiano. Cinética de processos fermentativos – Avaliação do perfil cinético dos cultivos por batelada, definição e cálculo dos parâmetros fermentativos (velocidades e rendimentos), modelos cinéticos para o crescimento celular (MONOD) e formação de produtos (GADEN). Formas de condução dos processos fermentativos - Principais configurações dos biorreatores, características e equacionamento dos processos descontínuos, contínuos e descontínuos- alimentado. Esterilização em processos fermentativos – Aspectos gerais sobre esterilização e desinfecção em bioprocessos, esterilização de meios de cultura, esterilização do ar, cinética da destruição térmica dos microrganismos, cálculo do tempo de esterilização por processo descontínuo e cálculo do tempo de esterilização por processo contínuo."
$styleC.Copy()
$ws.Cells.Item(17, 3).PasteSpecial(-4122)
$ws.Cells.Item(17, 3).Value = "Características do material biológico – Tipos de microrganismos, morfologia e estrutura celular, nutrição e crescimento microbiano. Cinética de processos fermentativos – Avaliação do perfil cinético dos cultivos por batelada, definição e cálculo dos parâmetros fermentativos (velocidades e rendimentos), modelos cinéticos para o crescimento celular (MONOD) e formação de produtos (GADEN). Formas de condução dos processos fermentativos - Principais configurações dos biorreatores, características e equacionamento dos processos descontínuos, contínuos e descontínuos- alimentado. Esterilização em processos fermentativos – Aspectos gerais sobre esterilização e desinfecção em bioprocessos, esterilização de meios de cultura, esterilização do ar, cinética da destruição térmica dos microrganismos, cálculo do tempo de esterilização por processo descontínuo e cálculo do tempo de esterilização por processo contínuo."
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Cells.Item(18, 1).Value = "Syllabus:"
$ws.Cells.Item(18, 2).Value = "Characteristics of biological material - Types of microorganisms, cell structure and morphology, nutrition and microbial growth. Kinetics of fermentative processes - Evaluation of kinetic profile of batch cultures, definition and calculation of fermentative parameters (rates and yields), kinetic models for cell growth (MONOD) and product formation (GADEN). Operation modes of fermentative processes. Major configurations of bioreactors, characteristics and mathematical equations for batch, fed-batch and continuous operations. Sterilization in fermentation process – general aspects on sterilization and disinfection in bioprocess, methods for medium and air sterilization, kinetics of thermal death of microorganisms, calculation of sterilization time for batch and continuous process."
$ws.Cells.Item(18, 3).Value = "Characteristics of biological material - Types of microorganisms, cell structure and morphology, nutrition and microbial growth. Kinetics of fermentative processes - Evaluation of kinetic profile of batch cultures, definition and calculation of fermentative parameters (rates and yields), kinetic models for cell growth (MONOD) and product formation (GADEN). Operation modes of fermentative processes. Major configurations of bioreactors, characteristics and mathematical equations for batch, fed-batch and continuous operations. Sterilization in fermentation process – general aspects on sterilization and disinfection in bioprocess, methods for medium and air sterilization, kinetics of thermal death of microorganisms, calculation of sterilization time for batch and continuous process."
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Cells.Item(19, 1).Value = "Avaliação:"
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).ClearContents()
$ws.Rows.Item(19).EntireRow.AutoFit()

# Row 20
$ws.Cells.Item(20, 1).Value = "Método:"
$ws.Cells.Item(20, 2).Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."
$ws.Cells.Item(20, 3).Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."

# Row 21
$ws.Cells.Item(21, 1).Value = "Critério:"
$ws.Cells.Item(21, 2).Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Cells.Item(21, 3).Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Cells.Item(22, 1).Value = "Norma de recuperação:"
$styleB.Copy()
$ws.Cells.Item(22, 2).PasteSpecial(-4122)
$ws.Cells.Item(22, 2).Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$styleC.Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4122)
$ws.Cells.Item(22, 3).Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$styleA.Copy()
$ws.Cells.Item(23, 1).PasteSpecial(-4122)
$ws.Cells.Item(23, 1).Value = "Bibliografia:"
$ws.Cells.Item(23, 2).Value = "1.Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – Engenharia Bioquímica, vol. 2, São Paulo: Edgard Blücher, 2001. 2. Borzani, W.; Schmidell, W.; Lima, U. A.; Aquarone, E. Biotecnologia Industrial. Fundamentos Vol. 1. São Paulo: Ed. Edgard Blücher, 2001. 3. Pauline M. Doran ed. Bioprocess Engineering Principles (Second Edition), Elsevier Ltd. 2013."
$ws.Cells.Item(23, 3).Value = "1.Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – Engenharia Bioquímica, vol. 2, São Paulo: Edgard Blücher, 2001. 2. Borzani, W.; Schmidell, W.; Lima, U. A.; Aquarone, E. Biotecnologia Industrial. Fundamentos Vol. 1. São Paulo: Ed. Edgard Blücher, 2001. 3. Pauline M. Doran ed. Bioprocess Engineering Principles (Second Edition), Elsevier Ltd. 2013."
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$styleA.Copy()
$ws.Cells.Item(24, 1).PasteSpecial(-4122)
$ws.Cells.Item(24, 1).Value = "Requisitos:"

# Row 25
$styleB.Copy()
$ws.Cells.Item(25, 2).PasteSpecial(-4122)
$ws.Cells.Item(25, 2).Value = "LOT2004 -  Bioquímica  (Requisito fraco)
"
$styleC.Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4122)
$ws.Cells.Item(25, 3).Value = "LOT2004 -  Bioquímica  (Requisito fraco)
"
$ws.Rows.Item(25).RowHeight = 30

# --- Column definitions: A alone gets width 30.71/style1; B width 60.71/style2 (unchanged) ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375

$ws.Application.CutCopyMode = $false
